# Update cryptocurrency price and 1h-volume-change figures in the
# worksheet to reflect the latest scrape (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$d = $ws.Cells.Item(2, 4)
$d.NumberFormat = '@'
$d.Value = '67.367.92'
$ws.Cells.Item(2, 5).Value = '  +0.31%  '
$d = $ws.Cells.Item(3, 4)
$d.NumberFormat = '@'
$d.Value = '2.633.77'
$ws.Cells.Item(3, 5).Value = '  +0.65%  '
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$d = $ws.Cells.Item(5, 4)
$d.NumberFormat = '@'
$d.Value = '602.66'
$ws.Cells.Item(5, 5).Value = '  +1.67%  '
$d = $ws.Cells.Item(6, 4)
$d.NumberFormat = '@'
$d.Value = '153.24'
$ws.Cells.Item(6, 5).Value = '  -0.58%  '
$ws.Cells.Item(7, 5).Value = '  +0.01%  '
$ws.Cells.Item(8, 5).Value = '  +3.36%  '
$d = $ws.Cells.Item(9, 4)
$d.NumberFormat = '@'
$d.Value = '2.631.02'
$ws.Cells.Item(9, 5).Value = '  +0.61%  '
$d = $ws.Cells.Item(10, 4)
$d.NumberFormat = '@'
$d.Value = '0.123'
$ws.Cells.Item(10, 5).Value = '  +3.75%  '
$ws.Cells.Item(11, 5).Value = '  +0.49%  '
$d = $ws.Cells.Item(12, 4)
$d.NumberFormat = '@'
$d.Value = '5.20'
$ws.Cells.Item(12, 5).Value = '  -0.35%  '
$d = $ws.Cells.Item(13, 4)
$d.NumberFormat = '@'
$d.Value = '0.352'
$ws.Cells.Item(13, 5).Value = '  -1.05%  '
$d = $ws.Cells.Item(14, 4)
$d.NumberFormat = '@'
$d.Value = '27.75'
$ws.Cells.Item(14, 5).Value = '  -0.11%  '
$d = $ws.Cells.Item(15, 4)
$d.NumberFormat = '@'
$d.Value = '3.111.91'
$ws.Cells.Item(15, 5).Value = '  +0.68%  '
$d = $ws.Cells.Item(16, 4)
$d.NumberFormat = '@'
$d.Value = '0.0000182'
$ws.Cells.Item(16, 5).Value = '  -0.77%  '
$d = $ws.Cells.Item(17, 4)
$d.NumberFormat = '@'
$d.Value = '67.381.23'
$ws.Cells.Item(17, 5).Value = '  +0.34%  '
$d = $ws.Cells.Item(18, 4)
$d.NumberFormat = '@'
$d.Value = '2.632.03'
$ws.Cells.Item(18, 5).Value = '  +0.66%  '
$d = $ws.Cells.Item(19, 4)
$d.NumberFormat = '@'
$d.Value = '11.23'
$ws.Cells.Item(19, 5).Value = '  -0.52%  '
$d = $ws.Cells.Item(20, 4)
$d.NumberFormat = '@'
$d.Value = '363.73'
$ws.Cells.Item(20, 5).Value = '  +0.73%  '
$d = $ws.Cells.Item(21, 4)
$d.NumberFormat = '@'
$d.Value = '7.56'
$ws.Cells.Item(21, 5).Value = '  -4.19%  '
$ws.Cells.Item(22, 5).Value = '  -0.50%  '
$ws.Cells.Item(23, 5).Value = '  +4.72%  '
$d = $ws.Cells.Item(24, 4)
$d.NumberFormat = '@'
$d.Value = '0.999'
$ws.Cells.Item(24, 5).Value = '  -0.12%  '
$d = $ws.Cells.Item(25, 4)
$d.NumberFormat = '@'
$d.Value = '10.17'
$ws.Cells.Item(25, 5).Value = '  -1.12%  '
$d = $ws.Cells.Item(26, 4)
$d.NumberFormat = '@'
$d.Value = '66.30'
$ws.Cells.Item(26, 5).Value = '  -7.30%  '
$d = $ws.Cells.Item(28, 4)
$d.NumberFormat = '@'
$d.Value = '0.0000103'
$ws.Cells.Item(28, 5).Value = '  -0.55%  '
$ws.Cells.Item(29, 5).Value = '  +0.50%  '
$d = $ws.Cells.Item(30, 4)
$d.NumberFormat = '@'
$d.Value = '576.99'
$ws.Cells.Item(30, 5).Value = '  -7.51%  '
$d = $ws.Cells.Item(31, 4)
$d.NumberFormat = '@'
$d.Value = '1.40'
$ws.Cells.Item(31, 5).Value = '  -3.97%  '
$d = $ws.Cells.Item(32, 4)
$d.NumberFormat = '@'
$d.Value = '7.87'
$ws.Cells.Item(32, 5).Value = '  -1.72%  '
$ws.Cells.Item(33, 5).Value = '  -0.18%  '
$d = $ws.Cells.Item(34, 4)
$d.NumberFormat = '@'
$d.Value = '0.128'
$ws.Cells.Item(34, 5).Value = '  -3.95%  '
$d = $ws.Cells.Item(35, 4)
$d.NumberFormat = '@'
$d.Value = '0.998'
$ws.Cells.Item(35, 5).Value = '  -0.05%  '
$ws.Cells.Item(36, 5).Value = '  -2.33%  '
$ws.Cells.Item(37, 5).Value = '  -1.37%  '
$d = $ws.Cells.Item(38, 4)
$d.NumberFormat = '@'
$d.Value = '158.05'
$ws.Cells.Item(38, 5).Value = '  +2.68%  '
$d = $ws.Cells.Item(39, 4)
$d.NumberFormat = '@'
$d.Value = '19.43'
$ws.Cells.Item(39, 5).Value = '  -0.05%  '
$ws.Cells.Item(40, 5).Value = '  +0.03%  '
$ws.Cells.Item(41, 5).Value = '  -3.75%  '
$ws.Cells.Item(42, 5).Value = '  -0.95%  '
$d = $ws.Cells.Item(43, 4)
$d.NumberFormat = '@'
$d.Value = '2.61'
$ws.Cells.Item(43, 5).Value = '  +0.53%  '
$d = $ws.Cells.Item(44, 4)
$d.NumberFormat = '@'
$d.Value = '41.22'
$ws.Cells.Item(44, 5).Value = '  -0.36%  '
$ws.Cells.Item(45, 5).Value = '  -0.03%  '
$d = $ws.Cells.Item(46, 4)
$d.NumberFormat = '@'
$d.Value = '16.35'
$ws.Cells.Item(46, 5).Value = '  -1.01%  '
$d = $ws.Cells.Item(47, 4)
$d.NumberFormat = '@'
$d.Value = '155.93'
$ws.Cells.Item(47, 5).Value = '  +0.10%  '
$d = $ws.Cells.Item(48, 4)
$d.NumberFormat = '@'
$d.Value = '0.0₆0291'
$ws.Cells.Item(48, 5).Value = '  -2.85%  '
$d = $ws.Cells.Item(49, 4)
$d.NumberFormat = '@'
$d.Value = '3.73'
$ws.Cells.Item(49, 5).Value = '  -1.09%  '
$ws.Cells.Item(50, 5).Value = '  +0.09%  '
$d = $ws.Cells.Item(51, 4)
$d.NumberFormat = '@'
$d.Value = '20.77'
$ws.Cells.Item(51, 5).Value = '  -1.39%  '
